$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string text cells (rows 5-8, columns B & C)
$ws.Range("B5").Value = "MT mem RD BW test, dura= 3.012559, GB/sec= 18.216385"
$ws.Range("C5").Value = "MT mem RD BW test, dura= 3.012559, GB/sec= 18.216385"
$ws.Range("B6").Value = "MT L3 RD BW test, dura= 3.000009, GB/sec= 65.738517"
$ws.Range("C6").Value = "MT L3 RD BW test, dura= 3.000009, GB/sec= 65.738517"
$ws.Range("B7").Value = "MT L2 RD BW test, dura= 3.000003, GB/sec= 133.187461"
$ws.Range("C7").Value = "MT L2 RD BW test, dura= 3.000003, GB/sec= 133.187461"
$ws.Range("B8").Value = "MT spin test, dura= 3.000000, Gops/sec= 3.857277"
$ws.Range("C8").Value = "MT spin test, dura= 3.000000, Gops/sec= 3.857277"

# Update numeric data cells (rows 5-8)
$ws.Range("D5").Value = 828.367259
$ws.Range("E5").Value = 831.374318
$ws.Range("G5").Value = 400
$ws.Range("H5").Value = 49.8
$ws.Range("I5").Value = 2.4
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 2.4
$ws.Range("L5").Value = 2.4
$ws.Range("M5").Value = 18.2
$ws.Range("N5").Value = 9.800000000000001
$ws.Range("P5").Value = 13.9
$ws.Range("R5").Value = 18242.7
$ws.Range("S5").Value = 18327.7
$ws.Range("T5").Value = 85
$ws.Range("U5").Value = 72.3
$ws.Range("V5").Value = 0.6
$ws.Range("W5").Value = 0.6
$ws.Range("X5").Value = 64.09999999999999
$ws.Range("Y5").Value = 63.9
$ws.Range("Z5").Value = 6.6
$ws.Range("AA5").Value = 6.6
$ws.Range("AB5").Value = 20.8
$ws.Range("AC5").Value = 20.8
$ws.Range("AF5").Value = 11.2
$ws.Range("AG5").Value = 11.1
$ws.Range("AH5").Value = 3.8
$ws.Range("AI5").Value = 3.8
$ws.Range("AJ5").Value = 0.6
$ws.Range("AK5").Value = 0.8
$ws.Range("AL5").Value = 3
$ws.Range("AM5").Value = 2.3
$ws.Range("AN5").Value = 0.6
$ws.Range("AO5").Value = 0.6
$ws.Range("AQ5").Value = 0.03
$ws.Range("AR5").Value = 0.35
$ws.Range("AS5").Value = 0.35
$ws.Range("AT5").Value = 0.03
$ws.Range("AU5").Value = 0.03
$ws.Range("AV5").Value = 88.2
$ws.Range("AW5").Value = 87.7
$ws.Range("AX5").Value = 88.40000000000001
$ws.Range("AY5").Value = 88.09999999999999
$ws.Range("AZ5").Value = 0.7
$ws.Range("BA5").Value = 1.1
$ws.Range("BB5").Value = 0.42
$ws.Range("BC5").Value = 0.43
$ws.Range("BD5").Value = 0.07000000000000001
$ws.Range("BE5").Value = 0.07000000000000001
$ws.Range("BF5").Value = 0.07000000000000001
$ws.Range("BG5").Value = 0.07000000000000001
$ws.Range("BH5").Value = 0.04
$ws.Range("BI5").Value = 0.04
$ws.Range("BJ5").Value = 0.04
$ws.Range("BK5").Value = 0.04
$ws.Range("BL5").Value = 0.01
$ws.Range("BN5").Value = 0.07000000000000001
$ws.Range("BO5").Value = 0.07000000000000001
$ws.Range("BP5").Value = 0.12
$ws.Range("BQ5").Value = 0.13
$ws.Range("BT5").Value = 0.16
$ws.Range("BU5").Value = 0.16
$ws.Range("BV5").Value = 0.17
$ws.Range("BW5").Value = 0.17
$ws.Range("BX5").Value = 0.09
$ws.Range("BY5").Value = 0.09
$ws.Range("BZ5").Value = 0.09
$ws.Range("CA5").Value = 0.09
$ws.Range("CC5").Value = 0.02
$ws.Range("CD5").Value = 0.17
$ws.Range("CE5").Value = 0.17
$ws.Range("CF5").Value = 0.3
$ws.Range("CG5").Value = 0.3
$ws.Range("CH5").Value = 0.01
$ws.Range("D6").Value = 831.380709
$ws.Range("E6").Value = 834.375218
$ws.Range("H6").Value = 56.8
$ws.Range("M6").Value = 22.6
$ws.Range("R6").Value = 519.1
$ws.Range("T6").Value = 5.8
$ws.Range("U6").Value = 0.2
$ws.Range("X6").Value = 29.3
$ws.Range("Y6").Value = 29.2
$ws.Range("AB6").Value = 43.4
$ws.Range("AC6").Value = 43.1
$ws.Range("AF6").Value = 25.8
$ws.Range("AG6").Value = 25.5
$ws.Range("AJ6").Value = 0.1
$ws.Range("AM6").Value = 0.7
$ws.Range("AQ6").Value = 0.03
$ws.Range("AR6").Value = 1.18
$ws.Range("AV6").Value = 66.8
$ws.Range("AW6").Value = 66.7
$ws.Range("AY6").Value = 66.09999999999999
$ws.Range("BA6").Value = 1.3
$ws.Range("BC6").Value = 1.23
$ws.Range("BO6").Value = 0.19
$ws.Range("CA6").Value = 0.28
$ws.Range("CB6").Value = 0.02
$ws.Range("CC6").Value = 0.01
$ws.Range("CD6").Value = 0.5
$ws.Range("CE6").Value = 0.5
$ws.Range("D7").Value = 834.388636
$ws.Range("E7").Value = 837.383139
$ws.Range("G7").Value = 399.7
$ws.Range("H7").Value = 61.5
$ws.Range("M7").Value = 23.1
$ws.Range("P7").Value = 16.7
$ws.Range("R7").Value = 493.5
$ws.Range("S7").Value = 498.6
$ws.Range("T7").Value = 5.1
$ws.Range("X7").Value = 0.2
$ws.Range("Y7").Value = 0.1
$ws.Range("Z7").Value = 1.5
$ws.Range("AA7").Value = 1
$ws.Range("AB7").Value = 41.2
$ws.Range("AC7").Value = 40.4
$ws.Range("AF7").Value = 17.6
$ws.Range("AG7").Value = 17.2
$ws.Range("AH7").Value = 25.8
$ws.Range("AI7").Value = 26
$ws.Range("AL7").Value = 0.1
$ws.Range("AN7").Value = 0.2
$ws.Range("AR7").Value = 2.37
$ws.Range("AS7").Value = 2.38
$ws.Range("AU7").Value = 0.01
$ws.Range("AV7").Value = 32
$ws.Range("AW7").Value = 31.6
$ws.Range("AX7").Value = 34.2
$ws.Range("AY7").Value = 34
$ws.Range("AZ7").Value = 1.7
$ws.Range("BA7").Value = 2.1
$ws.Range("BB7").Value = 2.47
$ws.Range("BC7").Value = 2.5
$ws.Range("BD7").Value = 0.38
$ws.Range("BE7").Value = 0.38
$ws.Range("BF7").Value = 0.39
$ws.Range("BG7").Value = 0.39
$ws.Range("BH7").Value = 0.21
$ws.Range("BI7").Value = 0.21
$ws.Range("BJ7").Value = 0.21
$ws.Range("BK7").Value = 0.21
$ws.Range("BN7").Value = 0.39
$ws.Range("BP7").Value = 0.64
$ws.Range("BQ7").Value = 0.64
$ws.Range("BT7").Value = 0.98
$ws.Range("BU7").Value = 0.99
$ws.Range("BV7").Value = 1
$ws.Range("BW7").Value = 1.01
$ws.Range("BX7").Value = 0.55
$ws.Range("BY7").Value = 0.55
$ws.Range("BZ7").Value = 0.55
$ws.Range("CA7").Value = 0.55
$ws.Range("CD7").Value = 1.01
$ws.Range("CE7").Value = 1.02
$ws.Range("CF7").Value = 1.65
$ws.Range("CG7").Value = 1.66
$ws.Range("CI7").Value = 0.02
$ws.Range("D8").Value = 837.392574
$ws.Range("E8").Value = 840.387074
$ws.Range("H8").Value = 65.5
$ws.Range("M8").Value = 23.1
$ws.Range("N8").Value = 13.1
$ws.Range("R8").Value = 482.8
$ws.Range("S8").Value = 487.3
$ws.Range("T8").Value = 4.5
$ws.Range("AB8").Value = 0.1
$ws.Range("AC8").Value = 0.1
$ws.Range("AF8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0.1
$ws.Range("AM8").Value = 0.1
$ws.Range("AP8").Value = 2.87
$ws.Range("AQ8").Value = 2.84
$ws.Range("AR8").Value = 0
$ws.Range("AS8").Value = 0
$ws.Range("AV8").Value = 11.4
$ws.Range("AW8").Value = 11.4
$ws.Range("AX8").Value = 8
$ws.Range("AY8").Value = 8.199999999999999
$ws.Range("AZ8").Value = 0.7
$ws.Range("BA8").Value = 0.8
$ws.Range("BB8").Value = 3.31
$ws.Range("BC8").Value = 3.31
$ws.Range("BD8").Value = 0.5600000000000001
$ws.Range("BO8").Value = 0.62
$ws.Range("BV8").Value = 1.53
$ws.Range("CD8").Value = 1.63
$ws.Range("CE8").Value = 1.61
$ws.Range("CF8").Value = 1.94
$ws.Range("CG8").Value = 1.93

Write-Host "done"